# working_hours.xlsx — add a missing work session (2014-02-22, 11:00-11:30)
# as a new row 17, pushing the trailing blank/summary rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17 (shifts the old row 17 onward down to 18+,
# and keeps the existing F2:F16 / G3:G16 shared formulas untouched
# since they sit entirely above the insertion point).
$ws.Rows(17).Insert()

# Fill in the new day's data.
$ws.Range("A17").Value = 2014
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 0.45833333333333331   # 11:00
$ws.Range("E17").Value = 0.47916666666666669   # 11:30

# The author's edit left A17 carrying the new time-only cell style
# (numFmtId 20, "h:mm") introduced for this row's cells - reproduce it
# by applying that number format after setting the value.
$ws.Range("A17").NumberFormat = "h:mm"

# Time-spent / hours-spent formulas, matching the pattern used by the
# rows above (F = minutes between start/end, G = F in hours).
$ws.Range("F17").Formula = "=(E17-D17)*24*60"
$ws.Range("G17").Formula = "=F17/60"

# Restore the active selection to the newly added row's F cell (the
# author's selection moved from F16 to F17 after the insert).
$ws.Range("F17").Select() | Out-Null
